$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A12's value (tiny correction in fractional seconds of the timestamp)
$ws.Range("A12").Value = 45874.41689232639

# Row 13 data
$ws.Range("A13").Value = 45874.43292002315
$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = 16.17
$ws.Range("E13").Value = 89.66
$ws.Range("F13").Value = 499.22
$ws.Range("G13").Value = 9.34
$ws.Range("H13").Value = "ESE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "10:23:24"

# Row 14 data
$ws.Range("A14").Value = 45874.44107065022
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = 19
$ws.Range("D14").Value = 17.22
$ws.Range("E14").Value = 85.18000000000001
$ws.Range("F14").Value = 335.74
$ws.Range("G14").Value = 9.16
$ws.Range("H14").Value = "ESE"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "10:35:08"

# Apply same style as A12 (date/time format) to the new A13/A14 cells
$ws.Range("A13:A14").NumberFormat = $ws.Range("A12").NumberFormat
